# Seperated and created new file test case for login with wrong credential
#
# Adds a second row to the "Login" sheet containing a new test email
# (a "wrong credential" login) with its own mailto hyperlink, and makes
# the Login sheet the active/selected sheet+cell (previously SignUpData
# was the selected tab).

$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item("Login")

# --- New data row for the "wrong credential" login test case ----------
$loginSheet.Range("A2").Value = "filty@yomail.com"
$loginSheet.Range("B2").Value = 12345678

# Excel always stamps newly-added hyperlinks with the built-in
# "Hyperlink" style, so re-apply the same look the row above already
# has (blue, non-underlined text) right after creating the link.
$loginSheet.Hyperlinks.Add($loginSheet.Range("A2"), "mailto:filty@yomail.com", $null, $null, "filty@yomail.com")
$loginSheet.Range("A1:B1").Copy()
$loginSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection / active sheet bookkeeping ------------------------------
$loginSheet.Activate()
$loginSheet.Range("B10").Select()
